$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- M1: "《二泉映月》" -> "《赛马》"
$ws.Range("M1").Value = "《赛马》"

# --- New column O: header + 10 names, copying the style/format of column M
$ws.Range("M1:M11").Copy()
$ws.Range("O1:O11").PasteSpecial(-4122)

$ws.Range("O1").Value = "《黄河颂》"
$ws.Range("O2").Value = "李昀杰"
$ws.Range("O3").Value = "魏伟涛"
$ws.Range("O4").Value = "王青粟"
$ws.Range("O5").Value = "赵子淳"
$ws.Range("O6").Value = "范菁宸"
$ws.Range("O7").Value = "王荣焕"
$ws.Range("O8").Value = "盛煜航"
$ws.Range("O9").Value = "李伟民"
$ws.Range("O10").Value = "黄思源"
$ws.Range("O11").Value = "冼永峰"

# Blank, but styled, cells O12:O16 (same style as rest of column)
$ws.Range("M12:M16").Copy()
$ws.Range("O12:O16").PasteSpecial(-4122)

# --- Selection moved to B7
$ws.Range("B7").Select()
